$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'248.86"
$ws.Range("G2").Value = "'10"

$ws.Range("D3").Value = "'22.67"
$ws.Range("G3").Value = "'10"

$ws.Range("D4").Value = "'5.406"
$ws.Range("G4").Value = "'10"

$ws.Range("D5").Value = "'0.05684"
$ws.Range("G5").Value = "'10"

$ws.Range("D6").Value = "'3.407"
$ws.Range("G6").Value = "'10"

$ws.Range("D7").Value = "'6.326"
$ws.Range("G7").Value = "'10"

$ws.Range("D8").Value = "'0.8049"
$ws.Range("G8").Value = "'10"

$ws.Range("D9").Value = "'0.9164"
$ws.Range("G9").Value = "'10"

$ws.Range("B10").Value = "One"
$ws.Range("C10").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D10").Value = "'0.01110"
$ws.Range("E10").Value = "9OneONE"
$ws.Range("G10").Value = "'10"

$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1401"
$ws.Range("E11").Value = "10WazirXWRX"
$ws.Range("G11").Value = "'10"

$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.07432"
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"
$ws.Range("G12").Value = "'10"

$ws.Range("B13").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C13").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D13").Value = "'0.03134"
$ws.Range("E13").Value = "12LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("G13").Value = "'10"

$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").Value = "'0.03033"
$ws.Range("E14").Value = "13BitrueCoinBTR"
$ws.Range("G14").Value = "'10"

$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Value = "'0.09374"
$ws.Range("E15").Value = "14BitMartTokenBMX"
$ws.Range("G15").Value = "'10"

$ws.Range("B16").Value = "MCDex"
$ws.Range("C16").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D16").Value = "'3.860"
$ws.Range("E16").Value = "15MCDexMCB"
$ws.Range("G16").Value = "'10"

$ws.Range("B17").Value = "BitForexToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D17").Value = "'0.001571"
$ws.Range("E17").Value = "16BitForexTokenBF"
$ws.Range("G17").Value = "'10"

$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D18").Value = "'0.04786"
$ws.Range("E18").Value = "17CoinExTokenCET"
$ws.Range("G18").Value = "'10"

$ws.Range("B19").Value = "UpBots"
$ws.Range("C19").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D19").Value = "'0.01827"
$ws.Range("E19").Value = "18UpBotsUBXTBestin24h"
$ws.Range("G19").Value = "'10"

$ws.Range("D20").Value = "'0.006452"
$ws.Range("G20").Value = "'10"

$ws.Range("D21").Value = "'0.004990"
$ws.Range("G21").Value = "'10"

$ws.Range("D22").Value = "'0.001007"
$ws.Range("G22").Value = "'10"

$ws.Range("G23").Value = "'10"

$ws.Range("D24").Value = "'3.700"
$ws.Range("G24").Value = "'10"

$ws.Range("D25").Value = "'2.201"
$ws.Range("G25").Value = "'10"

$ws.Range("D26").Value = "'0.3252"
$ws.Range("G26").Value = "'10"

$ws.Range("G27").Value = "'10"

$ws.Range("G28").Value = "'10"

$ws.Range("G29").Value = "'10"

$ws.Range("G30").Value = "'10"

$ws.Range("G31").Value = "'10"

$ws.Range("G32").Value = "'10"

$ws.Range("G33").Value = "'10"

$ws.Range("G34").Value = "'10"

$ws.Range("G35").Value = "'10"

$ws.Range("G36").Value = "'10"

$ws.Range("G37").Value = "'10"

$ws.Range("G38").Value = "'10"

$ws.Range("G39").Value = "'10"

$ws.Range("D40").Value = "'0.04005"
$ws.Range("G40").Value = "'10"

$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "'0.006810"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("G41").Value = "'10"

$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1069"
$ws.Range("E42").Value = "41BKEXTokenBKK"
$ws.Range("G42").Value = "'10"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "'0.002722"
$ws.Range("E43").Value = "42CEJICEJI"
$ws.Range("G43").Value = "'10"

$ws.Range("D44").Value = "'0.007966"
$ws.Range("G44").Value = "'10"

$ws.Range("D45").Value = "'0.00005794"
$ws.Range("G45").Value = "'10"

$ws.Range("D46").Value = "'0.00000000750"
$ws.Range("G46").Value = "'10"

$ws.Range("D47").Value = "'0.4990"
$ws.Range("G47").Value = "'10"

$ws.Range("D48").Value = "'0.2099"
$ws.Range("G48").Value = "'10"

$ws.Range("D49").Value = "'0.00002100"
$ws.Range("G49").Value = "'10"

$ws.Range("D50").Value = "'0.01010"
$ws.Range("G50").Value = "'10"

$ws.Range("G51").Value = "'10"
